$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.728.34'
$ws.Range('E2').Value = '  +3.30%  '
$ws.Range('D3').Value = '2.326.26'
$ws.Range('E3').Value = '  +1.32%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '521.42'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.86'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.89%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.538'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.39%  '
$ws.Range('D9').Value = '2.350.92'
$ws.Range('E9').Value = '  +1.41%  '
$ws.Range('E10').Value = '  +6.31%  '
$ws.Range('E11').Value = '  -0.73%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.24'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.16%  '
$ws.Range('E13').Value = '  +0.22%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.82'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.21%  '
$ws.Range('D15').Value = '2.746.62'
$ws.Range('E15').Value = '  +1.49%  '
$ws.Range('D16').Value = '56.783.77'
$ws.Range('E16').Value = '  +3.50%  '
$ws.Range('E17').Value = '  +2.37%  '
$ws.Range('D18').Value = '2.350.02'
$ws.Range('E18').Value = '  +3.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.50'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.11%  '
$ws.Range('E20').Value = '  +0.90%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '323.81'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.08%  '
$ws.Range('E22').Value = '  -1.42%  '
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.60'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.49%  '
$ws.Range('E25').Value = '  +8.54%  '
$ws.Range('E26').Value = '  +0.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.90'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.93%  '
$ws.Range('E28').Value = '  +13.47%  '
$ws.Range('D29').Value = '0.0₃0751'
$ws.Range('E29').Value = '  +5.76%  '
$ws.Range('E30').Value = '  +5.88%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '169.93'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.80%  '
$ws.Range('E32').Value = '  +0.34%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.31'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.26%  '
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.992'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.16%  '
$ws.Range('E36').Value = '  +1.17%  '
$ws.Range('E37').Value = '  +0.77%  '
$ws.Range('E38').Value = '  +3.37%  '
$ws.Range('E39').Value = '  +8.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '37.96'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.34%  '
$ws.Range('E41').Value = '  +0.26%  '
$ws.Range('E42').Value = '  +4.57%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '138.09'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.95%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.26'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.77%  '
$ws.Range('B45').Value = 'Bittensor'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '278.89'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0935'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.45%  '
$ws.Range('E47').Value = '  -0.09%  '
$ws.Range('E48').Value = '  +1.84%  '
$ws.Range('E49').Value = '  +3.45%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.83'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.41%  '
$ws.Range('E51').Value = '  +0.33%  '
